# heliaphen_experiments.xlsx — "corrected timezone error. added interpolated plot."
#
# 1. Header cell A1 "Experiment" -> "experiment" (lower-cased).
# 2. Selection moved from I3 back to A1.
# 3. Column widths bumped up (~+1 char of padding across the board).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the header text -------------------------------------------------
$ws.Range("A1").Value = "experiment"

# --- 2. Column widths --------------------------------------------------------
# ColumnWidth is in "characters"; the stored sheet width ends up
# ColumnWidth + 5/6. Values chosen so the saved width lands as close as
# possible to the target widths below:
#   A        -> 17.280612244898
#   B        -> 183.408163265306
#   C:D      -> 14.0408163265306
#   E:F      -> 14.7602040816327
#   G:H      -> 6.65816326530612
#   I        -> 19.4387755102041
#   J:AMK    -> 18.3571428571429
$ws.Columns.Item(1).ColumnWidth = 16.5
$ws.Columns.Item(2).ColumnWidth = 182.5
$ws.Range("C:D").ColumnWidth = 13.16666666666667
$ws.Range("E:F").ColumnWidth = 14
$ws.Range("G:H").ColumnWidth = 5.833333333333333
$ws.Columns.Item(9).ColumnWidth = 18.66666666666667
$ws.Range("J:AMK").ColumnWidth = 17.5

# --- 3. Selection -------------------------------------------------------------
[void]$ws.Range("A1").Select()
